$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.038386164815461
$ws.Range("D2").Value = 1.045851041287301
$ws.Range("E2").Value = 1.053256060915684
$ws.Range("F2").Value = 1.05892530057322
$ws.Range("I2").Value = 1.039873480429175
$ws.Range("J2").Value = 1.043483884721495
$ws.Range("K2").Value = 1.048618122783219
$ws.Range("L2").Value = 1.056002532289689
$ws.Range("M2").Value = 1.061656203973915
$ws.Range("N2").Value = 1.044965751527075
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.039199336134313
$ws.Range("D3").Value = 1.046489036361063
$ws.Range("E3").Value = 1.054132573153599
$ws.Range("F3").Value = 1.059770038747412
$ws.Range("I3").Value = 1.040047423652195
$ws.Range("J3").Value = 1.043942511211551
$ws.Range("K3").Value = 1.049068093421251
$ws.Range("L3").Value = 1.056691890253293
$ws.Range("M3").Value = 1.062314993400018
$ws.Range("N3").Value = 1.045425029319351
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.039726095241533
$ws.Range("D4").Value = 1.046902410562692
$ws.Range("E4").Value = 1.054701169430197
$ws.Range("F4").Value = 1.06031780566552
$ws.Range("I4").Value = 1.040159084803811
$ws.Range("J4").Value = 1.044239175539035
$ws.Range("K4").Value = 1.049359113278943
$ws.Range("L4").Value = 1.05713874290285
$ws.Range("M4").Value = 1.06274178966017
$ws.Range("N4").Value = 1.045722114944208
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.039947682615582
$ws.Range("D5").Value = 1.047076322650875
$ws.Range("E5").Value = 1.054940548728408
$ws.Range("F5").Value = 1.060548363773083
$ws.Range("I5").Value = 1.040205812964704
$ws.Range("J5").Value = 1.044363868577386
$ws.Range("K5").Value = 1.049481422909743
$ws.Range("L5").Value = 1.057326787721873
$ws.Range("M5").Value = 1.06292133656277
$ws.Range("N5").Value = 1.045846985060976
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.03998489614469
$ws.Range("D6").Value = 1.047105530801851
$ws.Range("E6").Value = 1.054980761530808
$ws.Range("F6").Value = 1.060587091673906
$ws.Range("I6").Value = 1.040213646248125
$ws.Range("J6").Value = 1.04438480360538
$ws.Range("K6").Value = 1.049501957137689
$ws.Range("L6").Value = 1.057358372249436
$ws.Range("M6").Value = 1.062951490369721
$ws.Range("N6").Value = 1.045867949819111
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.039729055563356
$ws.Range("D7").Value = 1.046904733877644
$ws.Range("E7").Value = 1.054704366688966
$ws.Range("F7").Value = 1.060320885308698
$ws.Range("I7").Value = 1.040159710030267
$ws.Range("J7").Value = 1.044240841791933
$ws.Range("K7").Value = 1.04936074772584
$ws.Range("L7").Value = 1.057141254830394
$ws.Range("M7").Value = 1.06274418829832
$ws.Range("N7").Value = 1.045723783563376
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.038660858156393
$ws.Range("D8").Value = 1.046066540014555
$ws.Range("E8").Value = 1.053551984193948
$ws.Range("F8").Value = 1.059210541424106
$ws.Range("I8").Value = 1.039932449575051
$ws.Range("J8").Value = 1.043638899158239
$ws.Range("K8").Value = 1.048770220925311
$ws.Range("L8").Value = 1.056235339082243
$ws.Range("M8").Value = 1.061878737228453
$ws.Range("N8").Value = 1.045120986102098
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.036783097324744
$ws.Range("D9").Value = 1.044593817186431
$ws.Range("E9").Value = 1.051532409447836
$ws.Range("F9").Value = 1.057262981986295
$ws.Range("I9").Value = 1.039525189017308
$ws.Range("J9").Value = 1.042577506988299
$ws.Range("K9").Value = 1.047728613342663
$ws.Range("L9").Value = 1.05464513732987
$ws.Range("M9").Value = 1.060357722400386
$ws.Range("N9").Value = 1.044058086633534
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.035534410236522
$ws.Range("D10").Value = 1.043614987919598
$ws.Range("E10").Value = 1.050193585366168
$ws.Range("F10").Value = 1.055970783412089
$ws.Range("I10").Value = 1.039249152976232
$ws.Range("J10").Value = 1.041869521955669
$ws.Range("K10").Value = 1.047033597862915
$ws.Range("L10").Value = 1.053589226810045
$ws.Range("M10").Value = 1.059346514018797
$ws.Range("N10").Value = 1.043349096180951
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.034994484255312
$ws.Range("D11").Value = 1.043191874930752
$ws.Range("E11").Value = 1.04961567720965
$ws.Range("F11").Value = 1.055412736543169
$ws.Range("I11").Value = 1.03912856096315
$ws.Range("J11").Value = 1.041562879939339
$ws.Range("K11").Value = 1.046732520181389
$ws.Range("L11").Value = 1.053133028383365
$ws.Range("M11").Value = 1.05890933489668
$ws.Range("N11").Value = 1.043042018697785
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.034794048169658
$ws.Range("D12").Value = 1.043034822921825
$ws.Range("E12").Value = 1.049401290544258
$ws.Range("F12").Value = 1.055205678012228
$ws.Range("I12").Value = 1.039083608178865
$ws.Range("J12").Value = 1.04144896872536
$ws.Range("K12").Value = 1.046620668111376
$ws.Range("L12").Value = 1.052963730559852
$ws.Range("M12").Value = 1.058747051067544
$ws.Range("N12").Value = 1.042927945716816
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.034837037092596
$ws.Range("D13").Value = 1.04306850605863
$ws.Range("E13").Value = 1.049447264773557
$ws.Range("F13").Value = 1.055250082559257
$ws.Range("I13").Value = 1.03909325791199
$ws.Range("J13").Value = 1.041473403531846
$ws.Range("K13").Value = 1.046644661572036
$ws.Range("L13").Value = 1.053000038495566
$ws.Range("M13").Value = 1.058781856768244
$ws.Range("N13").Value = 1.042952415223531
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.034977913760178
$ws.Range("D14").Value = 1.043178890687824
$ws.Range("E14").Value = 1.049597950336774
$ws.Range("F14").Value = 1.055395616420344
$ws.Range("I14").Value = 1.039124848402509
$ws.Range("J14").Value = 1.041553464214071
$ws.Range("K14").Value = 1.046723274822348
$ws.Range("L14").Value = 1.053119031009314
$ws.Range("M14").Value = 1.058895918334208
$ws.Range("N14").Value = 1.043032589601107
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.035064728008595
$ws.Range("D15").Value = 1.043246917063817
$ws.Range("E15").Value = 1.049690829073509
$ws.Range("F15").Value = 1.055485314500407
$ws.Range("I15").Value = 1.039144291226577
$ws.Range("J15").Value = 1.041602790860972
$ws.Range("K15").Value = 1.046771708640134
$ws.Range("L15").Value = 1.053192366758241
$ws.Range("M15").Value = 1.058966209255852
$ws.Range("N15").Value = 1.043081986297504
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.03557025958385
$ws.Range("D16").Value = 1.043643083999657
$ws.Range("E16").Value = 1.050231977564009
$ws.Range("F16").Value = 1.05600785057617
$ws.Range("I16").Value = 1.03925713385905
$ws.Range("J16").Value = 1.041889871203553
$ws.Range("K16").Value = 1.047053576742006
$ws.Range("L16").Value = 1.053619524783668
$ws.Range("M16").Value = 1.059375542644186
$ws.Range("N16").Value = 1.043369474327101
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.035887571878156
$ws.Range("D17").Value = 1.043891784867176
$ws.Range("E17").Value = 1.050571912116114
$ws.Range("F17").Value = 1.056336022291956
$ws.Range("I17").Value = 1.039327632000205
$ws.Range("J17").Value = 1.042069928647678
$ws.Range("K17").Value = 1.047230350914364
$ws.Range("L17").Value = 1.053887743400632
$ws.Range("M17").Value = 1.059632490116008
$ws.Range("N17").Value = 1.04354978747345
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.036072728348193
$ws.Range("D18").Value = 1.044036917853129
$ws.Range("E18").Value = 1.050770364839167
$ws.Range("F18").Value = 1.056527582315686
$ws.Range("I18").Value = 1.039368649416017
$ws.Range("J18").Value = 1.042174945320862
$ws.Range("K18").Value = 1.04733344753187
$ws.Range("L18").Value = 1.054044288872408
$ws.Range("M18").Value = 1.059782428801622
$ws.Range("N18").Value = 1.043654953282355
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.03613587434741
$ws.Range("D19").Value = 1.044086416264195
$ws.Range("E19").Value = 1.050838061652188
$ws.Range("F19").Value = 1.056592923529566
$ws.Range("I19").Value = 1.039382617810196
$ws.Range("J19").Value = 1.042210751932783
$ws.Range("K19").Value = 1.047368598608419
$ws.Range("L19").Value = 1.054097683419768
$ws.Range("M19").Value = 1.059833565089708
$ws.Range("N19").Value = 1.043690810743772
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.035853519649166
$ws.Range("D20").Value = 1.043865094366392
$ws.Range("E20").Value = 1.050535422276507
$ws.Range("F20").Value = 1.056300797765427
$ws.Range("I20").Value = 1.039320078861139
$ws.Range("J20").Value = 1.042050610984271
$ws.Range("K20").Value = 1.047211386029569
$ws.Range("L20").Value = 1.05385895592342
$ws.Range("M20").Value = 1.059604915289401
$ws.Range("N20").Value = 1.043530442376745
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.034936425864059
$ws.Range("D21").Value = 1.043146382087044
$ws.Range("E21").Value = 1.049553569610073
$ws.Range("F21").Value = 1.055352754100359
$ws.Range("I21").Value = 1.039115550188903
$ws.Range("J21").Value = 1.041529888624169
$ws.Range("K21").Value = 1.046700125677165
$ws.Range("L21").Value = 1.053083986399106
$ws.Range("M21").Value = 1.058862327162557
$ws.Range("N21").Value = 1.043008980531163
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.034360487278597
$ws.Range("D22").Value = 1.042695142238876
$ws.Range("E22").Value = 1.048937827036346
$ws.Range("F22").Value = 1.054757983758468
$ws.Range("I22").Value = 1.038986032186903
$ws.Range("J22").Value = 1.041202428442715
$ws.Range("K22").Value = 1.046378570123197
$ws.Range("L22").Value = 1.052597627506689
$ws.Range("M22").Value = 1.05839603486085
$ws.Range("N22").Value = 1.042681055318692
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.034665738694339
$ws.Range("D23").Value = 1.042934291449195
$ws.Range("E23").Value = 1.049264092712687
$ws.Range("F23").Value = 1.055073158647901
$ws.Range("I23").Value = 1.039054779357728
$ws.Range("J23").Value = 1.041376026670918
$ws.Range("K23").Value = 1.046549042399386
$ws.Range("L23").Value = 1.052855370144028
$ws.Range("M23").Value = 1.058643167576799
$ws.Range("N23").Value = 1.04285490007629
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.035868906158555
$ws.Range("D24").Value = 1.043877154438839
$ws.Range("E24").Value = 1.050551909924959
$ws.Range("F24").Value = 1.056316713771116
$ws.Range("I24").Value = 1.039323492117769
$ws.Range("J24").Value = 1.042059339829069
$ws.Range("K24").Value = 1.047219955484289
$ws.Range("L24").Value = 1.053871963442362
$ws.Range("M24").Value = 1.059617374964208
$ws.Range("N24").Value = 1.043539183617503
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.037267995123221
$ws.Range("D25").Value = 1.044974032877669
$ws.Range("E25").Value = 1.05205319408217
$ws.Range("F25").Value = 1.057765393125146
$ws.Range("I25").Value = 1.039631276983102
$ws.Range("J25").Value = 1.042851976601464
$ws.Range("K25").Value = 1.047998007173747
$ws.Range("L25").Value = 1.05505550523313
$ws.Range("M25").Value = 1.060750381293609
$ws.Range("N25").Value = 1.044332946025032

Write-Host "Applied 264 cell updates to vm_pu sheet"
